# Applies the "BBB spreadsheet cleanup results" update to the
# 11-2021 comparison deflators sheet: refreshed numeric values for a
# number of rows/quarters, plus clearing a few cells that no longer
# carry data for the 2021 Q2-Q4 columns on some rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated numeric values
$ws.Range("E5").Value = -0.2535
$ws.Range("F5").Value = -0.209
$ws.Range("G5").Value = -0.0857
$ws.Range("H5").Value = -0.1995
$ws.Range("I5").Value = -0.0735
$ws.Range("J5").Value = -0.1004
$ws.Range("K5").Value = -0.5243
$ws.Range("L5").Value = -0.2685
$ws.Range("M5").Value = -0.1031
$ws.Range("F7").Value = -0.0428
$ws.Range("G7").Value = 0.0535
$ws.Range("J7").Value = -0.1533
$ws.Range("K7").Value = -0.1994
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = -0.0295
$ws.Range("G9").Value = -0.0212
$ws.Range("H9").Value = -0.0212
$ws.Range("I9").Value = -0.0525
$ws.Range("J9").Value = -0.0401
$ws.Range("K9").Value = -0.0396
$ws.Range("L9").Value = -0.0248
$ws.Range("M9").Value = -0.0244
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.0027
$ws.Range("G10").Value = -0.005
$ws.Range("H10").Value = -0.0072
$ws.Range("I10").Value = -0.2309
$ws.Range("J10").Value = -0.1947
$ws.Range("K10").Value = -0.1809
$ws.Range("L10").Value = -0.1675
$ws.Range("M10").Value = -0.125
$ws.Range("F12").Value = -0.6014
$ws.Range("G12").Value = -0.6539
$ws.Range("H12").Value = -0.8608
$ws.Range("I12").Value = -0.6834
$ws.Range("J12").Value = 0.0193
$ws.Range("K12").Value = 0.0347
$ws.Range("L12").Value = 0.1981
$ws.Range("M12").Value = 0.0738
$ws.Range("E13").Value = -0.1236
$ws.Range("F13").Value = -0.1112
$ws.Range("G13").Value = -0.1281
$ws.Range("H13").Value = -0.1171
$ws.Range("I13").Value = -0.0491
$ws.Range("J13").Value = -0.0434
$ws.Range("K13").Value = -0.6121
$ws.Range("L13").Value = -0.6467
$ws.Range("M13").Value = -0.3035
$ws.Range("F15").Value = -2.0325
$ws.Range("G15").Value = -4.2436
$ws.Range("H15").Value = -2.9287
$ws.Range("I15").Value = -2.2538
$ws.Range("J15").Value = -1.8699
$ws.Range("K15").Value = -2.9846
$ws.Range("L15").Value = -1.2404
$ws.Range("M15").Value = -0.621
$ws.Range("E22").Value = -0.0269
$ws.Range("F22").Value = -0.0195
$ws.Range("G22").Value = -0.0203
$ws.Range("H22").Value = -0.0198
$ws.Range("I22").Value = -0.0196
$ws.Range("J22").Value = -0.0219
$ws.Range("K22").Value = -0.0229
$ws.Range("L22").Value = -0.0144
$ws.Range("M22").Value = -0.0129
$ws.Range("F23").Value = 0.0489
$ws.Range("G23").Value = -0.0288
$ws.Range("J23").Value = 0.0957
$ws.Range("K23").Value = 0.0718
$ws.Range("E26").Value = 159.2979
$ws.Range("F26").Value = 161.2187
$ws.Range("G26").Value = 163.1626
$ws.Range("H26").Value = 165.13
$ws.Range("I26").Value = 167.1211
$ws.Range("J26").Value = 169.1363
$ws.Range("K26").Value = 171.1757
$ws.Range("L26").Value = 173.2397
$ws.Range("M26").Value = 175.3286
$ws.Range("E32").Value = -0.0444
$ws.Range("F32").Value = -0.1066
$ws.Range("G32").Value = -0.1086
$ws.Range("H32").Value = -0.111
$ws.Range("I32").Value = -0.0333
$ws.Range("J32").Value = -0.1177
$ws.Range("K32").Value = -0.121
$ws.Range("L32").Value = -0.1244
$ws.Range("M32").Value = -0.0689
$ws.Range("F34").Value = -0.0708
$ws.Range("G34").Value = -0.0697
$ws.Range("H34").Value = 0.0022
$ws.Range("J34").Value = 0.0089
$ws.Range("K34").Value = -0.0206
$ws.Range("C36").Value = 0.0226
$ws.Range("D36").Value = -0.0448
$ws.Range("E36").Value = -0.0061
$ws.Range("F36").Value = -0.0265
$ws.Range("G36").Value = -0.0149
$ws.Range("H36").Value = 0.1092
$ws.Range("I36").Value = 0.0235
$ws.Range("J36").Value = 0.0038
$ws.Range("K36").Value = 0.0102
$ws.Range("L36").Value = -0.0035
$ws.Range("M36").Value = 0.007
$ws.Range("C37").Value = -0.1231
$ws.Range("D37").Value = -0.1382
$ws.Range("E37").Value = -0.2464
$ws.Range("F37").Value = -0.2147
$ws.Range("G37").Value = -0.1699
$ws.Range("H37").Value = -0.1295
$ws.Range("I37").Value = -0.0871
$ws.Range("J37").Value = -0.0582
$ws.Range("K37").Value = -0.0346
$ws.Range("L37").Value = -0.0222
$ws.Range("M37").Value = 0
$ws.Range("F39").Value = 0.3496
$ws.Range("G39").Value = -0.5157
$ws.Range("H39").Value = 0.0809
$ws.Range("I39").Value = 0.2205
$ws.Range("J39").Value = 0.216
$ws.Range("K39").Value = 0.2119
$ws.Range("L39").Value = 0.2079
$ws.Range("M39").Value = 0.068
$ws.Range("E40").Value = 0.0384
$ws.Range("F40").Value = 0.032
$ws.Range("G40").Value = 0.0294
$ws.Range("H40").Value = 0.0271
$ws.Range("I40").Value = -0.0111
$ws.Range("J40").Value = -0.006
$ws.Range("K40").Value = -0.0043
$ws.Range("L40").Value = -0.0027
$ws.Range("M40").Value = -0.0026
$ws.Range("F42").Value = -0.0844
$ws.Range("G42").Value = -0.9329
$ws.Range("H42").Value = -0.1691
$ws.Range("I42").Value = 0.0321
$ws.Range("J42").Value = 0.0097
$ws.Range("K42").Value = 0.013
$ws.Range("L42").Value = -0.0064
$ws.Range("M42").Value = -0.1061
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("F50").Value = 0.0714
$ws.Range("G50").Value = 0.0705
$ws.Range("I50").Value = -0.017
$ws.Range("J50").Value = -0.0989
$ws.Range("K50").Value = -0.1124
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0

# Cells that are cleared (no longer have a value) in the updated sheet
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("C39").ClearContents()
$ws.Range("D39").ClearContents()
$ws.Range("E39").ClearContents()
$ws.Range("C42").ClearContents()
$ws.Range("D42").ClearContents()
$ws.Range("E42").ClearContents()
